# Corrected excel sheets for application fix issues
#
# The "Edit Repayment Schedule" sheet had its approve/disburse rows
# (A8:B13 - "clickonapprove" ... "disburseloan") moved out to a brand
# new "Sheet1" tab (rows A1:B6 there). "Edit Repayment Schedule" keeps
# only its first 7 rows. Selections/active-tab are updated to match
# what Excel leaves behind after doing that move:
#   - "Edit Repayment Schedule": selection becomes the (now deleted)
#     former row block A8:XFD13, active cell A8, and it is no longer
#     the active tab.
#   - new "Sheet1": selection left at D5, not the active tab.
#   - "NewLoanInput" (first sheet) ends up active, with its selection
#     at B7.

$wb = $excel.ActiveWorkbook

$editSchedule = $wb.Worksheets.Item("Edit Repayment Schedule")

# New sheet goes right after "Edit Repayment Schedule" (last sheet today).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)

# Move (cut) the approve/disburse block to the new sheet.
$editSchedule.Range("A8:B13").Copy($newSheet.Range("A1:B6"))
$editSchedule.Range("A8:B13").EntireRow.Delete()

# Leftover selection on "Edit Repayment Schedule" after the row delete.
$editSchedule.Range("A8:XFD13").Select() | Out-Null

# Selection on the freshly created sheet.
$newSheet.Range("D5").Select() | Out-Null

# NewLoanInput becomes the active sheet/tab again, selection at B7.
$newLoanInput = $wb.Worksheets.Item("NewLoanInput")
$newLoanInput.Activate() | Out-Null
$newLoanInput.Range("B7").Select() | Out-Null
